$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swap (rows 25 and 26) ---
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

# --- Price column (D) updates: force text so values like "30.329.87" / "1.001" are not reinterpreted as numbers/dates ---
$priceUpdates = [ordered]@{
    "D2" = "30.329.87"
    "D3" = "1.871.60"
    "D4" = "1.001"
    "D5" = "236.20"
    "D6" = "1.001"
    "D8" = "0.2894"
    "D9" = "0.06641"
    "D10" = "21.71"
    "D11" = "0.08050"
    "D12" = "97.48"
    "D13" = "1.872.36"
    "D15" = "0.6888"
    "D16" = "272.28"
    "D17" = "30.323.48"
    "D18" = "14.18"
    "D19" = "0.000007800"
    "D20" = "1.001"
    "D21" = "2.117.13"
    "D23" = "1.002"
    "D25" = "9.359"
    "D26" = "168.27"
    "D27" = "19.01"
    "D28" = "1.960"
    "D30" = "0.09986"
    "D32" = "1.465"
    "D33" = "4.096"
    "D35" = "1.139"
    "D36" = "0.7033"
    "D37" = "2.715"
    "D38" = "0.01888"
    "D39" = "2.653"
    "D40" = "6.325"
    "D41" = "72.81"
    "D42" = "1.966"
    "D43" = "0.8444"
    "D44" = "0.4176"
    "D46" = "103.24"
    "D47" = "9.332"
    "D48" = "7.116"
    "D49" = "934.72"
    "D50" = "34.57"
    "D51" = "0.05681"
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}

# --- Volume(1h) column (E) updates: plain text percentages ---
$volumeUpdates = [ordered]@{
    "E2" = "  +0.19%  "
    "E3" = "  +0.46%  "
    "E4" = "  +0.05%  "
    "E5" = "  +0.17%  "
    "E6" = "  +0.06%  "
    "E7" = "  +0.89%  "
    "E8" = "  +1.90%  "
    "E9" = "  +1.72%  "
    "E10" = "  +0.04%  "
    "E11" = "  +1.39%  "
    "E12" = "  -0.05%  "
    "E13" = "  +0.41%  "
    "E14" = "  -0.28%  "
    "E16" = "  -2.61%  "
    "E17" = "  +0.17%  "
    "E18" = "  +5.83%  "
    "E19" = "  +6.73%  "
    "E21" = "  +0.13%  "
    "E22" = "  -1.20%  "
    "E23" = "  +0.10%  "
    "E24" = "  +0.82%  "
    "E25" = "  +1.94%  "
    "E26" = "  +0.49%  "
    "E27" = "  -0.39%  "
    "E28" = "  +1.50%  "
    "E29" = "  -0.96%  "
    "E30" = "  +2.47%  "
    "E31" = "  -0.14%  "
    "E32" = "  -0.82%  "
    "E33" = "  +0.73%  "
    "E35" = "  +0.81%  "
    "E37" = "  +0.19%  "
    "E38" = "  +1.08%  "
    "E39" = "  +2.77%  "
    "E40" = "  +0.25%  "
    "E41" = "  -2.96%  "
    "E42" = "  +0.26%  "
    "E43" = "  -0.75%  "
    "E44" = "  -0.17%  "
    "E45" = "  +0.08%  "
    "E46" = "  -0.30%  "
    "E47" = "  +0.29%  "
    "E48" = "  -0.97%  "
    "E49" = "  -3.19%  "
    "E50" = "  +1.35%  "
    "E51" = "  +0.58%  "
}
foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}
